$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 283, shifting existing rows 283:346 down to 284:347
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with the new record's data
$ws.Range("A283").Value2 = 6
$ws.Range("B283").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C283").Value2 = "Metropolitana"
$ws.Range("D283").Value2 = 44463
$ws.Range("E283").Value2 = 13
$ws.Range("F283").Value2 = 100112003
$ws.Range("G283").Value2 = "Ajo"
$ws.Range("H283").Value2 = "Chino"
$ws.Range("I283").Value2 = "Primera"
$ws.Range("J283").Value2 = 200
$ws.Range("K283").Value2 = 14000
$ws.Range("L283").Value2 = 15000
$ws.Range("M283").Value2 = 14500
$ws.Range("N283").Value2 = "$/caja 10 kilos"
$ws.Range("O283").Value2 = "China"
$ws.Range("P283").Value2 = 1450
$ws.Range("Q283").Value2 = 10
$ws.Range("R283").Value2 = "Hortaliza"
